# Refactor PPV Tools: add two columns ("Program" and "Lot") to the end
# of the "ppv" table on the PPV worksheet, formatted like the existing
# header/body cells, then leave the selection on the first new data
# cell (J2) -- matching the author's edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPV")
$lo = $ws.ListObjects.Item("ppv")

# Insert the "Program" column right after the existing last column (BinDesc).
$programCol = $lo.ListColumns.Add()
$programCol.Range.Cells.Item(1, 1).Value = "Program"

# Insert the "Lot" column right after "Program".
$lotCol = $lo.ListColumns.Add()
$lotCol.Range.Cells.Item(1, 1).Value = "Lot"

# Match the look of the rest of the table: copy the header format from
# the existing last header cell, and the body format from the existing
# last body cell, onto the two new columns.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("H2").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Move the selection to the first cell of the newly-added data, as the
# author's session left it.
$ws.Activate()
$ws.Range("J2").Select()
